$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("P2").Value = 266134856
$ws.Range("V2").Value = 973244716178

# ---- Row 3 ----
$ws.Range("P3").Value = 721090219
$ws.Range("V3").Value = 973222132639

# ---- Row 4 ----
$ws.Range("P4").Value = 640547479
$ws.Range("V4").Value = "6802B874E598"

# ---- Row 5 ----
$ws.Range("P5").Value = 824035814
# new blank numeric-styled cells R5:U5 (same style as V5 -> numFmt "0", right aligned, no border/font override)
$ws.Range("R5:U5").NumberFormat = "0"
$ws.Range("R5:U5").HorizontalAlignment = -4152
$ws.Range("V5").Value = 973037989617

# ---- Row 6 ----
$ws.Range("P6").Value = 824035814
$ws.Range("R6:U6").NumberFormat = "0"
$ws.Range("R6:U6").HorizontalAlignment = -4152
$ws.Range("V6").Value = 973161507193

# ---- Row 7: clear out almost everything, keep only styled-but-empty A7,B7,K7,N7,P7 ----
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7:J7").Clear()
$ws.Range("K7").ClearContents()
$ws.Range("L7:M7").Clear()
$ws.Range("N7").ClearContents()
$ws.Range("O7").Clear()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").Clear()
$ws.Range("V7").Clear()

# ---- Selection moved from D11 to D12 ----
[void]$ws.Range("D12").Select()
